$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.254.60'

$ws.Range("D3").Value = '2.346.17'
$ws.Range("E3").Value = '  -5.23%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''308.86'
$ws.Range("E5").Value = '  -3.71%  '

$ws.Range("D6").Value = '''85.96'
$ws.Range("E6").Value = '  -6.68%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '''0.488'
$ws.Range("E9").Value = '  -4.52%  '

$ws.Range("D10").Value = '''0.0823'

$ws.Range("D11").Value = '''30.50'
$ws.Range("E11").Value = '  -7.55%  '

$ws.Range("E12").Value = '  +0.41%  '

$ws.Range("D13").Value = '2.709.05'
$ws.Range("E13").Value = '  -5.20%  '

$ws.Range("D14").Value = '''6.47'
$ws.Range("E14").Value = '  -6.19%  '

$ws.Range("D15").Value = '''14.87'
$ws.Range("E15").Value = '  -3.93%  '

$ws.Range("D16").Value = '2.343.08'
$ws.Range("E16").Value = '  -5.74%  '

$ws.Range("D17").Value = '''0.757'
$ws.Range("E17").Value = '  -4.62%  '

$ws.Range("D18").Value = '40.223.11'
$ws.Range("E18").Value = '  -3.47%  '

$ws.Range("E19").Value = '  -3.51%  '

$ws.Range("D20").Value = '''6.12'
$ws.Range("E20").Value = '  -4.98%  '

$ws.Range("D21").Value = '''67.97'
$ws.Range("E21").Value = '  -4.63%  '

$ws.Range("D22").Value = '''10.85'
$ws.Range("E22").Value = '  -3.52%  '

$ws.Range("D23").Value = '''236.07'
$ws.Range("E23").Value = '  -1.44%  '

$ws.Range("E24").Value = '  -6.83%  '

$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").Value = '''1.81'
$ws.Range("E26").Value = '  -6.60%  '

$ws.Range("D27").Value = '''23.56'
$ws.Range("E27").Value = '  -5.71%  '

$ws.Range("E28").Value = '  -3.47%  '

$ws.Range("D29").Value = '''9.27'
$ws.Range("E29").Value = '  -4.90%  '

$ws.Range("D30").Value = '''35.06'
$ws.Range("E30").Value = '  -4.41%  '

$ws.Range("D31").Value = '''152.19'
$ws.Range("E31").Value = '  -3.15%  '

$ws.Range("E32").Value = '  -0.04%  '

$ws.Range("D33").Value = '''5.17'
$ws.Range("E33").Value = '  -4.95%  '

$ws.Range("D34").Value = '''0.0728'
$ws.Range("E34").Value = '  -5.36%  '

$ws.Range("D35").Value = '''2.44'
$ws.Range("E35").Value = '  -4.69%  '

$ws.Range("E36").Value = '  -2.16%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.101'
$ws.Range("E37").Value = '  -2.49%  '

$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").Value = '''15.96'
$ws.Range("E38").Value = '  -6.97%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''2.75'
$ws.Range("E39").Value = '  -4.24%  '

$ws.Range("D40").Value = '''1.72'
$ws.Range("E40").Value = '  -6.61%  '

$ws.Range("D41").Value = '''3.83'
$ws.Range("E41").Value = '  -4.27%  '

$ws.Range("E42").Value = '  -5.60%  '

$ws.Range("D43").Value = '1.956.54'
$ws.Range("E43").Value = '  -2.17%  '

$ws.Range("D44").Value = '''0.0269'
$ws.Range("E44").Value = '  -5.32%  '

$ws.Range("D45").Value = '''17.69'
$ws.Range("E45").Value = '  -4.97%  '

$ws.Range("D46").Value = '''9.34'
$ws.Range("E46").Value = '  -1.00%  '

$ws.Range("D47").Value = '''2.70'
$ws.Range("E47").Value = '  -8.93%  '

$ws.Range("D48").Value = '2.567.58'
$ws.Range("E48").Value = '  -6.21%  '

$ws.Range("D49").Value = '''93.42'
$ws.Range("E49").Value = '  -4.40%  '

$ws.Range("D50").Value = '''71.65'
$ws.Range("E50").Value = '  -5.93%  '

$ws.Range("D51").Value = '''50.91'
$ws.Range("E51").Value = '  -2.40%  '
